$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new review row (row 9), following the same layout/formatting as the
# previous row (row 8): copy its formatting down first, then set the values.
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)

$ws.Range("A9").Value = "com.singleton.strechy"
$ws.Range("B9").Value = "taxi game"
$ws.Range("C9").Value = "galiatia942@gmail.com"
$ws.Range("D9").Value = "syechimovitz@gmail.com"
$ws.Range("E9").Value = "27/5/2019 15:59"
$ws.Range("F9").Value = "this app car taxi game is a gift. I must recommend it to every kid and parent. Challenge yourself now!"

# Match the row height used by the other "taxi game" review rows (7, 8).
$ws.Rows.Item(9).RowHeight = 13.8

$null = $ws.Range("F9").Select()
